$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row58 = 58
$ws.Cells.Item($row58, 1).Value = 45929
$ws.Cells.Item($row58, 2).Value = "四方坪站"
$ws.Cells.Item($row58, 3).Value = 9170.6299999999992
$ws.Cells.Item($row58, 4).Value = 7445.21
$ws.Cells.Item($row58, 5).Value = 3207.95
$ws.Cells.Item($row58, 6).Value = 373

$row59 = 59
$ws.Cells.Item($row59, 1).Value = 45929
$ws.Cells.Item($row59, 2).Value = "高岭站"
$ws.Cells.Item($row59, 3).Value = 4957.18
$ws.Cells.Item($row59, 4).Value = 3947.78
$ws.Cells.Item($row59, 5).Value = 1241.19
$ws.Cells.Item($row59, 6).Value = 189
